$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (the issue URLs) before re-adding the new ones
$ws.Hyperlinks.Delete()

$ws.Range("E2").Value = "https://github.com/GIScience/ohsome-api/pull/6"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/GIScience/ohsome-api/pull/6")
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Value = "https://github.com/GIScience/oshdb/pull/258"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/GIScience/oshdb/pull/258")
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Value = "https://github.com/GIScience/ohsome-api/pull/33"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/GIScience/ohsome-api/pull/33")
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Value = "https://github.com/GIScience/ohsome-api/pull/27"
$ws.Hyperlinks.Add($ws.Range("E5"), "https://github.com/GIScience/ohsome-api/pull/27")
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Value = "https://github.com/GIScience/ohsome-api/pull/60"
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/GIScience/ohsome-api/pull/60")
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E7").Value = "https://github.com/GIScience/ohsome-api/commit/c54cd1bedacefd0eac5674fff932fddd2e5c2232"
$ws.Hyperlinks.Add($ws.Range("E7"), "https://github.com/GIScience/ohsome-api/commit/c54cd1bedacefd0eac5674fff932fddd2e5c2232")
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("E8").Value = "https://github.com/GIScience/ohsome-api/pull/151"
$ws.Hyperlinks.Add($ws.Range("E8"), "https://github.com/GIScience/ohsome-api/pull/151")
$ws.Range("E8").Style = "Hyperlink"
$ws.Range("E9").Value = "https://github.com/GIScience/ohsome-api/pull/112"
$ws.Hyperlinks.Add($ws.Range("E9"), "https://github.com/GIScience/ohsome-api/pull/112")
$ws.Range("E9").Style = "Hyperlink"
$ws.Range("E10").Value = "https://github.com/GIScience/ohsome-api/pull/130"
$ws.Hyperlinks.Add($ws.Range("E10"), "https://github.com/GIScience/ohsome-api/pull/130")
$ws.Range("E10").Style = "Hyperlink"
$ws.Range("E11").Value = "https://github.com/GIScience/ohsome-api/pull/144"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://github.com/GIScience/ohsome-api/pull/144")
$ws.Range("E11").Style = "Hyperlink"
$ws.Range("E12").Value = "https://github.com/GIScience/ohsome-api/commit/cda684d0aa7fb748ebe4205610f94b3961de4797"
$ws.Hyperlinks.Add($ws.Range("E12"), "https://github.com/GIScience/ohsome-api/commit/cda684d0aa7fb748ebe4205610f94b3961de4797")
$ws.Range("E12").Style = "Hyperlink"
$ws.Range("E13").Value = "https://github.com/GIScience/ohsome-api/pull/215"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://github.com/GIScience/ohsome-api/pull/215")
$ws.Range("E13").Style = "Hyperlink"
$ws.Range("E14").Value = "https://github.com/GIScience/ohsome-api/pull/294"
$ws.Hyperlinks.Add($ws.Range("E14"), "https://github.com/GIScience/ohsome-api/pull/294")
$ws.Range("E14").Style = "Hyperlink"
$ws.Range("E15").Value = "https://github.com/GIScience/ohsome-api/commit/44c01dccf285b2b1b494f0764f8f0f8bc987c362"
$ws.Hyperlinks.Add($ws.Range("E15"), "https://github.com/GIScience/ohsome-api/commit/44c01dccf285b2b1b494f0764f8f0f8bc987c362")
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("E16").Value = "https://github.com/GIScience/ohsome-api/pull/314"
$ws.Hyperlinks.Add($ws.Range("E16"), "https://github.com/GIScience/ohsome-api/pull/314")
$ws.Range("E16").Style = "Hyperlink"
$ws.Range("E17").Value = "https://github.com/GIScience/ohsome-api/pull/319"
$ws.Hyperlinks.Add($ws.Range("E17"), "https://github.com/GIScience/ohsome-api/pull/319")
$ws.Range("E17").Style = "Hyperlink"
$ws.Range("E18").Value = "https://github.com/confluentinc/kafka-rest/pull/144"
$ws.Hyperlinks.Add($ws.Range("E18"), "https://github.com/confluentinc/kafka-rest/pull/144")
$ws.Range("E18").Style = "Hyperlink"
$ws.Range("E19").Value = "https://github.com/confluentinc/kafka-rest/pull/222"
$ws.Hyperlinks.Add($ws.Range("E19"), "https://github.com/confluentinc/kafka-rest/pull/222")
$ws.Range("E19").Style = "Hyperlink"
$ws.Range("E20").Value = "https://github.com/confluentinc/kafka-rest/commit/e9c7bb73fb99519d4c38f824dd927687a6426466"
$ws.Hyperlinks.Add($ws.Range("E20"), "https://github.com/confluentinc/kafka-rest/commit/e9c7bb73fb99519d4c38f824dd927687a6426466")
$ws.Range("E20").Style = "Hyperlink"
$ws.Range("E21").Value = "https://github.com/confluentinc/kafka-rest/pull/67"
$ws.Hyperlinks.Add($ws.Range("E21"), "https://github.com/confluentinc/kafka-rest/pull/67")
$ws.Range("E21").Style = "Hyperlink"
$ws.Range("E22").Value = "https://github.com/fabioformosa/quartz-manager/commit/04da4556b14ca1d9f8a47406e6efc034a5e2ffd0"
$ws.Hyperlinks.Add($ws.Range("E22"), "https://github.com/fabioformosa/quartz-manager/commit/04da4556b14ca1d9f8a47406e6efc034a5e2ffd0")
$ws.Range("E22").Style = "Hyperlink"
$ws.Range("E23").Value = "https://github.com/senzing-garage/senzing-api-server/pull/391/commits"
$ws.Hyperlinks.Add($ws.Range("E23"), "https://github.com/senzing-garage/senzing-api-server/pull/391/commits")
$ws.Range("E23").Style = "Hyperlink"

# Widen column E to fit the longer commit/pull-request URLs
$ws.Columns.Item(5).ColumnWidth = 44.16666666666667

# Restore the active selection recorded in the workbook
$ws.Range("I8").Select()
